$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-number-looking Price cells to stay as text (matching the
# original inline-string cells) instead of being auto-converted to numbers,
# which would silently drop significant trailing/leading zeros.
$textCells = @("D4", "D5", "D6", "D8", "D11", "D15", "D17", "D18", "D20", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D32", "D34", "D39", "D40", "D41", "D42", "D44", "D45", "D46", "D47", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values scraped for this run.
$ws.Range("D2").Value = '70.018.39'
$ws.Range("E2").Value = '  +1.14%  '
$ws.Range("D3").Value = '3.539.68'
$ws.Range("E3").Value = '  +1.56%  '
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").Value = '605.54'
$ws.Range("E5").Value = '  -1.00%  '
$ws.Range("D6").Value = '197.43'
$ws.Range("E6").Value = '  +6.43%  '
$ws.Range("E7").Value = '  +0.40%  '
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("E9").Value = '  -6.62%  '
$ws.Range("E10").Value = '  +0.21%  '
$ws.Range("D11").Value = '53.89'
$ws.Range("E11").Value = '  +1.59%  '
$ws.Range("E12").Value = '  -1.01%  '
$ws.Range("E13").Value = '  +0.17%  '
$ws.Range("D14").Value = '4.089.43'
$ws.Range("E14").Value = '  +0.96%  '
$ws.Range("D15").Value = '599.00'
$ws.Range("E15").Value = '  -0.44%  '
$ws.Range("D16").Value = '70.143.01'
$ws.Range("E16").Value = '  +1.03%  '
$ws.Range("D17").Value = '19.14'
$ws.Range("E17").Value = '  +1.84%  '
$ws.Range("D18").Value = '12.73'
$ws.Range("E18").Value = '  +1.12%  '
$ws.Range("D19").Value = '3.540.40'
$ws.Range("E19").Value = '  +0.66%  '
$ws.Range("D20").Value = '0.121'
$ws.Range("E20").Value = '  +1.27%  '
$ws.Range("E21").Value = '  +0.88%  '
$ws.Range("D22").Value = '18.53'
$ws.Range("E22").Value = '  +7.53%  '
$ws.Range("E23").Value = '  +4.88%  '
$ws.Range("B24").Value = 'PancakeSwap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D24").Value = '4.64'
$ws.Range("E24").Value = '  -0.31%  '
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").Value = '101.68'
$ws.Range("E25").Value = '  -3.09%  '
$ws.Range("D26").Value = '3.19'
$ws.Range("E26").Value = '  +6.04%  '
$ws.Range("D27").Value = '10.91'
$ws.Range("E27").Value = '  -0.13%  '
$ws.Range("D28").Value = '9.64'
$ws.Range("E28").Value = '  -2.79%  '
$ws.Range("D29").Value = '33.42'
$ws.Range("E29").Value = '  -0.15%  '
$ws.Range("E30").Value = '  +11.67%  '
$ws.Range("E31").Value = '  +2.13%  '
$ws.Range("D32").Value = '12.52'
$ws.Range("E32").Value = '  +1.09%  '
$ws.Range("E33").Value = '  -1.13%  '
$ws.Range("D34").Value = '63.15'
$ws.Range("E34").Value = '  -0.24%  '
$ws.Range("D35").Value = '0.0₃0860'
$ws.Range("E35").Value = '  +11.51%  '
$ws.Range("D36").Value = '3.722.12'
$ws.Range("E36").Value = '  +4.11%  '
$ws.Range("E37").Value = '  +0.10%  '
$ws.Range("E38").Value = '  -2.81%  '
$ws.Range("D39").Value = '3.65'
$ws.Range("E39").Value = '  +2.01%  '
$ws.Range("D40").Value = '0.394'
$ws.Range("E40").Value = '  -0.42%  '
$ws.Range("D41").Value = '36.71'
$ws.Range("E41").Value = '  +0.22%  '
$ws.Range("D42").Value = '490.90'
$ws.Range("E42").Value = '  -5.62%  '
$ws.Range("E43").Value = '  -3.20%  '
$ws.Range("D44").Value = '0.0455'
$ws.Range("E44").Value = '  -0.95%  '
$ws.Range("B45").Value = 'ThetaToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D45").Value = '2.85'
$ws.Range("E45").Value = '  -3.47%  '
$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").Value = '0.140'
$ws.Range("E46").Value = '  -2.02%  '
$ws.Range("D47").Value = '3.30'
$ws.Range("E47").Value = '  -0.43%  '
$ws.Range("E48").Value = '  +0.21%  '
$ws.Range("D49").Value = '8.59'
$ws.Range("E49").Value = '  -2.63%  '
$ws.Range("D50").Value = '0.000251'
$ws.Range("E50").Value = '  +3.16%  '
$ws.Range("D51").Value = '130.87'
$ws.Range("E51").Value = '  -0.34%  '
